$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Asset Id","Asset Status","External Id","Product Id","Provider Id","Provider Name",
  "Marketplace Id","Marketplace Name","Contract Id","Contract Name","Reseller Id",
  "Reseller External Id","Reseller Name","Created At","Customer Id","Customer External Id",
  "Customer Name","Seamless Move","Discount Group","Action","Renewal Date","PurchaseType",
  "Adobe Customer Id","Adobe VIP Number","Adobe User Email","Currency","Cost","Reseller Cost",
  "MSRP","Seats","USD Cost","USD Reseller Cost","USD MSRP"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange = $ws.Range("A1:AG1")
$headerRange.Interior.ThemeColor = 1
$headerRange.Interior.TintAndShade = -0.14999847407452621
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.HorizontalAlignment = -4108

$ws.Range("A1:AG1").AutoFilter()

$ws.Select()
$ws.Range("H3").Select()

$wb.Save()
